$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.917.85"
$ws.Range("E2").Value = "'  +0.00%  "
$ws.Range("D3").Value = "'1.629.46"
$ws.Range("E3").Value = "'  -0.54%  "
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'211.89"
$ws.Range("E6").Value = "'  -0.15%  "
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("D8").Value = "'23.37"
$ws.Range("E8").Value = "'  -1.06%  "
$ws.Range("E9").Value = "'  -1.82%  "
$ws.Range("E11").Value = "'  +0.51%  "
$ws.Range("D12").Value = "'1.860.94"
$ws.Range("E12").Value = "'  -0.55%  "
$ws.Range("D13").Value = "'1.627.83"
$ws.Range("E13").Value = "'  -0.71%  "
$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = "'  -1.52%  "
$ws.Range("D15").Value = "'0.562"
$ws.Range("E15").Value = "'  -2.43%  "
$ws.Range("D16").Value = "'65.62"
$ws.Range("E16").Value = "'  -0.53%  "
$ws.Range("D17").Value = "'27.914.87"
$ws.Range("E17").Value = "'  +0.02%  "
$ws.Range("D18").Value = "'230.37"
$ws.Range("E18").Value = "'  -0.68%  "
$ws.Range("D19").Value = "'0.0₃0723"
$ws.Range("E19").Value = "'  -0.20%  "
$ws.Range("E20").Value = "'  +0.32%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "'  -0.04%  "
$ws.Range("E22").Value = "'  -0.33%  "
$ws.Range("D23").Value = "'10.26"
$ws.Range("E23").Value = "'  -4.84%  "
$ws.Range("E24").Value = "'  -1.23%  "
$ws.Range("D25").Value = "'154.85"
$ws.Range("E25").Value = "'  +2.06%  "
$ws.Range("E26").Value = "'  +0.17%  "
$ws.Range("E27").Value = "'  -0.15%  "
$ws.Range("D28").Value = "'15.54"
$ws.Range("E28").Value = "'  -1.28%  "
$ws.Range("E30").Value = "'  -0.78%  "
$ws.Range("E31").Value = "'  -0.36%  "
$ws.Range("E32").Value = "'  +2.32%  "
$ws.Range("E33").Value = "'  -0.74%  "
$ws.Range("D34").Value = "'1.401.17"
$ws.Range("E34").Value = "'  -0.63%  "
$ws.Range("E35").Value = "'  +0.16%  "
$ws.Range("E36").Value = "'  +11.02%  "
$ws.Range("E37").Value = "'  -0.18%  "
$ws.Range("E38").Value = "'  +2.01%  "
$ws.Range("E39").Value = "'  +0.10%  "
$ws.Range("E40").Value = "'  -3.11%  "
$ws.Range("E41").Value = "'  -0.03%  "
$ws.Range("E42").Value = "'  -0.03%  "
$ws.Range("D43").Value = "'1.85"
$ws.Range("E43").Value = "'  +0.26%  "
$ws.Range("E44").Value = "'  +0.10%  "
$ws.Range("D46").Value = "'2.20"
$ws.Range("E46").Value = "'  -0.36%  "
$ws.Range("D47").Value = "'1.771.03"
$ws.Range("E47").Value = "'  -0.49%  "
$ws.Range("D48").Value = "'88.01"
$ws.Range("E48").Value = "'  -0.34%  "
$ws.Range("B49").Value = "'Algorand"
$ws.Range("C49").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.101"
$ws.Range("E49").Value = "'  +0.45%  "
$ws.Range("B50").Value = "'Cronos"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0505"
$ws.Range("E50").Value = "'  -0.20%  "
$ws.Range("B51").Value = "'EnergySwap"
$ws.Range("C51").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.54"
$ws.Range("E51").Value = "'  -1.43%  "
